# This workbook tracks weekly wholesale-market price records for "Coliflor"
# (cauliflower). The edit adds one new week's worth of records (2 rows, for
# "Primera" and "Segunda" quality) right after the most recent existing
# record (row 1107), shifting every subsequent row down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1108-1109, pushing old rows 1108.. down to 1110..
$ws.Range("A1108:A1109").EntireRow.Insert()

# New row 1108: "Primera" quality record for date serial 45041 (2023-04-25)
$ws.Range("A1108").Value2 = 6
$ws.Range("B1108").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1108").Value2 = "Metropolitana"
$ws.Range("D1108").Value2 = 45041
$ws.Range("E1108").Value2 = 13
$ws.Range("F1108").Value2 = 100112008
$ws.Range("G1108").Value2 = "Coliflor"
$ws.Range("H1108").Value2 = "Sin especificar"
$ws.Range("I1108").Value2 = "Primera"
$ws.Range("J1108").Value2 = 7800
$ws.Range("K1108").Value2 = 1000
$ws.Range("L1108").Value2 = 1200
$ws.Range("M1108").Value2 = 1108
$ws.Range("N1108").Value2 = '$/unidad'
$ws.Range("O1108").Value2 = "Región Metropolitana"
$ws.Range("P1108").Value2 = 1108
$ws.Range("Q1108").Value2 = 1
$ws.Range("R1108").Value2 = "Hortaliza"

# New row 1109: "Segunda" quality record for the same date serial 45041
$ws.Range("A1109").Value2 = 6
$ws.Range("B1109").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1109").Value2 = "Metropolitana"
$ws.Range("D1109").Value2 = 45041
$ws.Range("E1109").Value2 = 13
$ws.Range("F1109").Value2 = 100112008
$ws.Range("G1109").Value2 = "Coliflor"
$ws.Range("H1109").Value2 = "Sin especificar"
$ws.Range("I1109").Value2 = "Segunda"
$ws.Range("J1109").Value2 = 2600
$ws.Range("K1109").Value2 = 800
$ws.Range("L1109").Value2 = 800
$ws.Range("M1109").Value2 = 800
$ws.Range("N1109").Value2 = '$/unidad'
$ws.Range("O1109").Value2 = "Región Metropolitana"
$ws.Range("P1109").Value2 = 800
$ws.Range("Q1109").Value2 = 1
$ws.Range("R1109").Value2 = "Hortaliza"
